# Insert a new data row at row 26, shifting the existing rows 26-66 down to 27-67.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26 (Excel copies formatting from the row above by default)
$ws.Rows.Item(26).Insert()

# Populate the new row 26 with the new data point
$ws.Cells.Item(26, 1).Value = 4
$ws.Cells.Item(26, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(26, 3).Value = "Los Lagos"
$ws.Cells.Item(26, 4).Value = 44645
$ws.Cells.Item(26, 5).Value = 10
$ws.Cells.Item(26, 6).Value = 100112031
$ws.Cells.Item(26, 7).Value = "Poroto verde"
$ws.Cells.Item(26, 8).Value = "Magnum"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 50
$ws.Cells.Item(26, 11).Value = 30000
$ws.Cells.Item(26, 12).Value = 30000
$ws.Cells.Item(26, 13).Value = 30000
$ws.Cells.Item(26, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(26, 15).Value = "Región Metropolitana"
$ws.Cells.Item(26, 16).Value = 1200
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"
